$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{
        "G" = 6.206015333333333
        "H" = 18.618046
        "I" = 0.0150172404156507
        "J" = 0.0150172404156507
        "K" = 2
        "L" = 0.6666666666666666
        "M" = 0.051093
        "N" = 0.153279
        "O" = 0.01450579975525089
        "P" = 0.01450579975525089
        "Q" = 0.317083941426
        "R" = 2.853755472834
        "S" = 0.0002178370823458896
        "T" = 0.0002178370823458897
    }
    3 = @{
        "G" = 6.206015333333333
        "H" = 18.618046
        "I" = 0.0150172404156507
        "J" = 0.0150172404156507
        "O" = 0.2313022967634575
        "P" = 0.2313022967634575
        "Q" = 5.056063447456221
        "R" = 45.504571027106
        "S" = 0.003473522199189025
        "T" = 0.003473522199189026
    }
    4 = @{
        "G" = 6.206015333333333
        "H" = 18.618046
        "I" = 0.0150172404156507
        "J" = 0.0150172404156507
        "M" = 2.656449666666667
        "N" = 7.969348999999999
        "O" = 0.7541919034812916
        "P" = 0.7541919034812917
        "Q" = 16.48596736356156
        "R" = 148.373706272054
        "S" = 0.01132588113411578
        "T" = 0.01132588113411578
    }
    5 = @{
        "I" = 0.9317452840597572
        "J" = 0.9317452840597571
        "K" = 2
        "L" = 0.6666666666666666
        "M" = 0.051093
        "N" = 0.153279
        "O" = 0.01450579975525089
        "P" = 0.01450579975525089
        "Q" = 19.673485866741
        "R" = 177.061372800669
        "S" = 0.01351571051347019
        "T" = 0.0135157105134702
    }
    6 = @{
        "I" = 0.9317452840597572
        "J" = 0.9317452840597571
        "O" = 0.2313022967634575
        "P" = 0.2313022967634575
        "S" = 0.2155148242015419
        "T" = 0.2155148242015419
    }
    7 = @{
        "I" = 0.9317452840597572
        "J" = 0.9317452840597571
        "M" = 2.656449666666667
        "N" = 7.969348999999999
        "O" = 0.7541919034812916
        "P" = 0.7541919034812917
        "Q" = 1022.87250646616
        "R" = 9205.85255819544
        "S" = 0.7027147493447451
        "T" = 0.7027147493447451
    }
    8 = @{
        "G" = 22.00088566666667
        "H" = 66.002657
        "I" = 0.05323747552459213
        "J" = 0.05323747552459213
        "K" = 2
        "L" = 0.6666666666666666
        "M" = 0.051093
        "N" = 0.153279
        "O" = 0.01450579975525089
        "P" = 0.01450579975525089
        "Q" = 1.124091251367
        "R" = 10.116821262303
        "S" = 0.0007722521594348035
        "T" = 0.0007722521594348036
    }
    9 = @{
        "G" = 22.00088566666667
        "H" = 66.002657
        "I" = 0.05323747552459213
        "J" = 0.05323747552459213
        "O" = 0.2313022967634575
        "P" = 0.2313022967634575
        "Q" = 17.92420222254744
        "R" = 161.317820002927
        "S" = 0.01231395036272651
        "T" = 0.01231395036272651
    }
    10 = @{
        "G" = 22.00088566666667
        "H" = 66.002657
        "I" = 0.05323747552459213
        "J" = 0.05323747552459213
        "M" = 2.656449666666667
        "N" = 7.969348999999999
        "O" = 0.7541919034812916
        "P" = 0.7541919034812917
        "Q" = 58.4442453955881
        "R" = 525.9982085602929
        "S" = 0.04015127300243081
        "T" = 0.04015127300243081
    }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
